$d = $word.ActiveDocument

# --- Step 1: remove the stray "_GoBack" bookmark that wraps the Banking
#     System UML picture paragraph (both bookmarkStart and bookmarkEnd). ---
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Step 2: rewrite the "Q2" paragraph. The single sentence becomes two
#     paragraphs, broken into several runs (mirroring Word's grammar-check
#     proofErr spans), and the "_GoBack" bookmark re-appears near the end
#     of the second paragraph. ---
$q2 = $d.Content
[void]$q2.Find.Execute(
    "Q2.We can create abstract class and after inheritance to child class.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Q2.We </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>can</w:t></w:r>
            <w:r><w:t xml:space="preserve">not </w:t></w:r>
            <w:r><w:t xml:space="preserve"> create</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> the object of </w:t></w:r>
            <w:r><w:t xml:space="preserve"> abstract class </w:t></w:r>
            <w:r><w:t>.</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t xml:space="preserve">But </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>Using  inheritance</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> in child class</w:t></w:r>
            <w:r><w:t xml:space="preserve"> to abstract class</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t xml:space="preserve"> we can create object.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$q2.InsertXML($xml)
